$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D6").Value = "Peace lily, Mint pot, Ficus, Mini red rose pot"
$ws.Range("D7").Value = "Aloe Vera, Mini palm, Mini pine tree"
$ws.Range("D8").Value = "Aloe Vera, Purple kalanchoe, Spiky boi"
$ws.Range("D3").Value = "Peace lily, Ficus, Aloe Vera, Mini palm"

$ws.Columns.Item(4).ColumnWidth = 37
